# Rename the inline logo pictures embedded in the document's headers and
# footers so the BTEC logo (a .jpg) swaps its "image1.jpg"/"image2.jpg"
# label and the Pearson logo (a .png) swaps its "image2.png"/"image1.png"
# label, in both the default and first-page header/footer.
#
# Section.Headers/Footers collection index 1 == "default" header/footer,
# index 2 == "first page" header/footer.
#
# NOTE: re-fetching the shape through its own .Range.InlineShapes before
# renaming keeps the handle fresh (this runtime stales out an
# InlineShapes(1) handle obtained straight off a HeaderFooter.Range once a
# prior edit has touched the package), so every shape is addressed the
# same, reliable way.

# --- Header (default) : BTec_Logo-Orange  image1.jpg -> image2.jpg
$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$raw = $sec.Headers.Item(1).Range.InlineShapes.Item(1)
$shp = $raw.Range.InlineShapes.Item(1)
$shp.Name = "image2.jpg"

# --- Header (first page) : BTec_Logo-Orange  image1.jpg -> image2.jpg
$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$raw = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$shp = $raw.Range.InlineShapes.Item(1)
$shp.Name = "image2.jpg"

# --- Footer (default) : PearsonLogo  image2.png -> image1.png
$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$raw = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$shp = $raw.Range.InlineShapes.Item(1)
$shp.Name = "image1.png"

# --- Footer (first page) : PearsonLogo  image2.png -> image1.png
$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$raw = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$shp = $raw.Range.InlineShapes.Item(1)
$shp.Name = "image1.png"
